$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the frequency-table values in rows 2-5 (columns B-X) with the
# recomputed fractions from the updated run, per "Updated run for publication".
$ws.Cells.Item(2, 2).Value = 0.00147710487444609
$ws.Cells.Item(2, 3).Value = 0.99753815854259
$ws.Cells.Item(2, 5).Value = 0.996553421959626
$ws.Cells.Item(2, 6).Value = 0.00147710487444609
$ws.Cells.Item(2, 7).Value = 0.000492368291482029
$ws.Cells.Item(2, 8).Value = 0.000492368291482029
$ws.Cells.Item(2, 11).Value = 0.789758739537174
$ws.Cells.Item(2, 12).Value = 0.178237321516494
$ws.Cells.Item(2, 14).Value = 0.998030526834072
$ws.Cells.Item(2, 15).Value = 0.0334810438207779
$ws.Cells.Item(2, 17).Value = 0.978335795174791
$ws.Cells.Item(2, 18).Value = 0.946824224519941
$ws.Cells.Item(2, 19).Value = 0.00147710487444609
$ws.Cells.Item(2, 20).Value = 0.0147710487444609
$ws.Cells.Item(2, 21).Value = 0.00295420974889217
$ws.Cells.Item(2, 23).Value = 0.00787789266371246
$ws.Cells.Item(2, 24).Value = 0.000492368291482029
$ws.Cells.Item(3, 3).Value = 0.000492368291482029
$ws.Cells.Item(3, 5).Value = 0.000492368291482029
$ws.Cells.Item(3, 6).Value = 0.000492368291482029
$ws.Cells.Item(3, 7).Value = 0.99753815854259
$ws.Cells.Item(3, 8).Value = 0.000492368291482029
$ws.Cells.Item(3, 9).Value = 0.00246184145741014
$ws.Cells.Item(3, 11).Value = 0.000492368291482029
$ws.Cells.Item(3, 12).Value = 0.000492368291482029
$ws.Cells.Item(3, 13).Value = 0.000984736582964057
$ws.Cells.Item(3, 15).Value = 0.00886262924667651
$ws.Cells.Item(3, 18).Value = 0.0507139340226489
$ws.Cells.Item(3, 20).Value = 0.953717380600689
$ws.Cells.Item(3, 21).Value = 0.000984736582964057
$ws.Cells.Item(3, 22).Value = 0.0206794682422452
$ws.Cells.Item(3, 23).Value = 0.0733628754308223
$ws.Cells.Item(4, 2).Value = 0.998030526834072
$ws.Cells.Item(4, 3).Value = 0.00196947316592811
$ws.Cells.Item(4, 5).Value = 0.00196947316592811
$ws.Cells.Item(4, 6).Value = 0.998030526834072
$ws.Cells.Item(4, 7).Value = 0.000984736582964057
$ws.Cells.Item(4, 8).Value = 0.999015263417036
$ws.Cells.Item(4, 11).Value = 0.168882323978336
$ws.Cells.Item(4, 12).Value = 0.821270310192024
$ws.Cells.Item(4, 14).Value = 0.00196947316592811
$ws.Cells.Item(4, 17).Value = 0.0108321024126046
$ws.Cells.Item(4, 18).Value = 0.000984736582964057
$ws.Cells.Item(4, 19).Value = 0.996553421959626
$ws.Cells.Item(4, 20).Value = 0.000984736582964057
$ws.Cells.Item(4, 21).Value = 0.99507631708518
$ws.Cells.Item(4, 22).Value = 0.978335795174791
$ws.Cells.Item(4, 23).Value = 0.891678975873954
$ws.Cells.Item(4, 24).Value = 0.999015263417036
$ws.Cells.Item(5, 2).Value = 0.000492368291482029
$ws.Cells.Item(5, 5).Value = 0.000984736582964057
$ws.Cells.Item(5, 7).Value = 0.000984736582964057
$ws.Cells.Item(5, 9).Value = 0.99753815854259
$ws.Cells.Item(5, 11).Value = 0.0408665681930084
$ws.Cells.Item(5, 13).Value = 0.999015263417036
$ws.Cells.Item(5, 15).Value = 0.957656326932546
$ws.Cells.Item(5, 17).Value = 0.0108321024126046
$ws.Cells.Item(5, 18).Value = 0.00147710487444609
$ws.Cells.Item(5, 19).Value = 0.00196947316592811
$ws.Cells.Item(5, 20).Value = 0.0305268340718858
$ws.Cells.Item(5, 21).Value = 0.000984736582964057
$ws.Cells.Item(5, 22).Value = 0.000984736582964057
$ws.Cells.Item(5, 23).Value = 0.0270802560315116
$ws.Cells.Item(5, 24).Value = 0.000492368291482029
